# feat: add scheduler example and alarm script
#
# Appends three new translation rows (VENTILATION pilot screen + schedule
# manager) to the "Translations" sheet, right after the existing
# "ELEMENTS IN MAINTENANCE STATE" entry at row 184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 185-187 -------------------------------------------------
# Column A first, then column B, for each pair, so new shared-string
# entries are appended in the same order as the source edit.

$ws.Range("A185").Value = "PILOTAGE VENTILATION"
$ws.Range("A186").Value = "PLANNING DE VENTILATION"

$ws.Range("B185").Value = "VENTILATION CONTROL"
$ws.Range("B186").Value = "VENTILATION SCHEDULE"

$ws.Range("A187").Value = "GESTIONNAIRE DE PLANNING"
$ws.Range("B187").Value = "SCHEDULE MANAGER"

# --- Selection / dimension ---------------------------------------------
# Excel recalculates dimension/used-range automatically; nudge the active
# selection to the next empty row like the real edit did.
[void]$ws.Range("B188").Select()
